$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-30
$dataI = @(6,2,5,6,6,5,7,6,5,7,7,7,7,6,4,1,5,3,4,1,8,4,1,7,7,4,7,3,7)
$dataJ = @(6,3,6,7,7,6,8,6,6,8,7,7,8,6,7,4,6,6,7,4,8,4,4,9,8,5,9,4,7)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
